$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 3 more rows of the same data (rows 5, 6, 7), matching rows 2-4
$name = "Videnov"
$location = "Sofia, Tsarigradsko, 15"

for ($r = 5; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $location
}
